$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Bug fix: numerical plantcodes were being stored/read as text (string),
# causing a type error downstream. Store them as real numbers instead.
$ws.Range("B6").Value = 4552
$ws.Range("B7").Value = 5740

# Move the active selection down to B8, matching the post-edit cursor
# position left by the author after editing B7.
$ws.Range("B8").Select()
